$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 4).Value = "Lietotājs veiksmīgi pieslēdzas"
